# Weekly refresh: the data rows (2-22) on Sheet1 get re-shuffled — each
# row's editable columns (Fecha, Variedad, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion,
# Origen, Precio $/Kg) are reassigned to a different row position while the
# fixed descriptive columns (Mercado ID/Mercado/Region/Codreg/Tipo/Producto.../
# Kg por unidad) stay put (they are identical for every row anyway).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a group, by 1-based column index:
# D=4 (Fecha), K=11 (Variedad), L=12 (Calidad), M=13 (Volumen),
# N=14 (Precio minimo), O=15 (Precio maximo), P=16 (Precio promedio ponderado),
# Q=17 (Unidad de comercializacion), R=18 (Origen), S=19 (Precio $/Kg)
$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19)

$firstRow = 2
$lastRow = 22

# Snapshot every row's values for the moving columns BEFORE any writes,
# since the mapping below moves data between rows.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# destinationRow -> sourceRow : destination row ends up with the source
# row's original (pre-edit) values for the moving columns.
$mapping = @{
    2  = 14
    3  = 15
    4  = 21
    5  = 20
    6  = 9
    7  = 19
    8  = 17
    9  = 8
    10 = 22
    11 = 18
    12 = 13
    13 = 2
    14 = 5
    15 = 7
    16 = 10
    17 = 4
    18 = 11
    19 = 6
    20 = 16
    21 = 3
    22 = 12
}

foreach ($destRow in ($mapping.Keys | Sort-Object)) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
